# Update the NATMI LR-pairs sheet with refreshed TPM-derived specificity scores.
# The sending-cluster pairing shifts (FAPs/MuSCs -> ECs/FAPs block) and every
# downstream derived-specificity column is recomputed against the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..J ligand stats, K..P receptor stats, Q..T edge weights/specificities.
$rows = @(
    @{ Row=2; A="ECs";  B="Rspo3"; C="Fzd8"; D="ECs";           E=1; F="0.3333333333333333";  G="0.06729733333333333"; H="0.201892";          I="0.01373511018321553"; J="0.01373511018321553"; K=3; L=1;                    M="3.390429";          N="10.171287";          O="0.173121426386348";  P="0.173121426386348";  Q="0.228166830556";     R="2.053501475004";    S="0.002377841866491927"; T="0.002377841866491927" },
    @{ Row=3; A="ECs";  B="Rspo3"; C="Fzd8"; D="FAPs";          E=1; F="0.3333333333333333";  G="0.06729733333333333"; H="0.201892";          I="0.01373511018321553"; J="0.01373511018321553"; K=3; L=1;                    M="11.625787";         N="34.877361";          O="0.5936336753560868"; P="0.5936336753560868"; Q="0.7823844630013334"; R="7.041460167012";    S="0.008153623939483051"; T="0.008153623939483051" },
    @{ Row=4; A="ECs";  B="Rspo3"; C="Fzd8"; D="MuSCs";         E=1; F="0.3333333333333333";  G="0.06729733333333333"; H="0.201892";          I="0.01373511018321553"; J="0.01373511018321553"; K=3; L=1;                    M="4.546141666666667"; N="13.638425";          O="0.2321342018628743"; P="0.2321342018628743"; Q="0.3059432111222223"; R="2.7534889001";      S="0.003188388839879374"; T="0.003188388839879374" },
    @{ Row=5; A="ECs";  B="Rspo3"; C="Fzd8"; D="Resolving-Mac"; E=1; F="0.3333333333333333";  G="0.06729733333333333"; H="0.201892";          I="0.01373511018321553"; J="0.01373511018321553"; K=2; L="0.6666666666666666"; M="0.021752";          N="0.06525600000000001"; O="0.001110696394691009"; P="0.001110696394691009"; Q="0.001463851594666667"; R="0.013174664352";   S="1.525553736118126E-05"; T="1.525553736118126E-05" },
    @{ Row=6; A="FAPs"; B="Rspo3"; C="Fzd8"; D="ECs";           E=3; F=1;                    G="4.83236";             H="14.49708";          I="0.9862648898167845"; J="0.9862648898167844"; K=3; L=1;                    M="3.390429";          N="10.171287";          O="0.173121426386348";  P="0.173121426386348";  Q="16.38377348244";     R="147.45396134196";   S="0.1707435845198561";   T="0.1707435845198561" },
    @{ Row=7; A="FAPs"; B="Rspo3"; C="Fzd8"; D="FAPs";          E=3; F=1;                    G="4.83236";             H="14.49708";          I="0.9862648898167845"; J="0.9862648898167844"; K=3; L=1;                    M="11.625787";         N="34.877361";          O="0.5936336753560868"; P="0.5936336753560868"; Q="56.17998806732001";  R="505.61989260588";   S="0.5854800514166038";   T="0.5854800514166038" },
    @{ Row=8; A="FAPs"; B="Rspo3"; C="Fzd8"; D="MuSCs";         E=3; F=1;                    G="4.83236";             H="14.49708";          I="0.9862648898167845"; J="0.9862648898167844"; K=3; L=1;                    M="4.546141666666667"; N="13.638425";          O="0.2321342018628743"; P="0.2321342018628743"; Q="21.96859314433334";  R="197.717338299";     S="0.2289458130229949";   T="0.2289458130229949" },
    @{ Row=9; A="FAPs"; B="Rspo3"; C="Fzd8"; D="Resolving-Mac"; E=3; F=1;                    G="4.83236";             H="14.49708";          I="0.9862648898167845"; J="0.9862648898167844"; K=2; L="0.6666666666666666"; M="0.021752";          N="0.06525600000000001"; O="0.001110696394691009"; P="0.001110696394691009"; Q="0.10511349472";      R="0.9460214524800001"; S="0.001095440857329828"; T="0.001095440857329828" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = [double]$r.F
    $ws.Range("G$n").Value = [double]$r.G
    $ws.Range("H$n").Value = [double]$r.H
    $ws.Range("I$n").Value = [double]$r.I
    $ws.Range("J$n").Value = [double]$r.J
    $ws.Range("K$n").Value = $r.K
    $ws.Range("L$n").Value = [double]$r.L
    $ws.Range("M$n").Value = [double]$r.M
    $ws.Range("N$n").Value = [double]$r.N
    $ws.Range("O$n").Value = [double]$r.O
    $ws.Range("P$n").Value = [double]$r.P
    $ws.Range("Q$n").Value = [double]$r.Q
    $ws.Range("R$n").Value = [double]$r.R
    $ws.Range("S$n").Value = [double]$r.S
    $ws.Range("T$n").Value = [double]$r.T
}
